$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows: OT values that moved from "Pendiente ADM" to an actual ticket id ---
$ws.Range("E7").Value = "ICD30377283"
$ws.Range("E8").Value = "ICD30377406"
$ws.Range("E12").Value = "ICD30377539"

# --- Append two new case rows (15 and 16) ---

# Row 15
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "5843"
$ws.Range("A15").ClearFormats()
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "8/13/2025"
$ws.Range("B15").ClearFormats()
$ws.Range("C15").Value = "CAMPANA 2545"
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = "Pendiente ADM"
$ws.Range("F15").Value = "Optical Power"
$ws.Range("G15").Value = "Pendiente"
$ws.Range("H15").Value = "Cable en panza"
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = '{"direccionesNormalizadas": [{"altura": 2545, "cod_calle": 3039, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.492468", "y": "-34.607061"}, "direccion": "CAMPANA 2545, CABA", "nombre_calle": "CAMPANA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K15").Value = -58.492468
$ws.Range("L15").Value = -34.607061
$ws.Range("M15").Value = "Devoto"
$ws.Range("N15").Value = "Capital Norte"

# Row 16
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "5860"
$ws.Range("A16").ClearFormats()
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "8/13/2025"
$ws.Range("B16").ClearFormats()
$ws.Range("C16").Value = "ARGERICH 45"
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = "Pendiente ADM"
$ws.Range("F16").Value = "Optical Power"
$ws.Range("G16").Value = "Pendiente"
$ws.Range("H16").Value = "cables a baja altura"
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = '{"direccionesNormalizadas": [{"altura": 45, "cod_calle": 1110, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.471436", "y": "-34.630493"}, "direccion": "ARGERICH 45, CABA", "nombre_calle": "ARGERICH", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K16").Value = -58.471436
$ws.Range("L16").Value = -34.630493
$ws.Range("M16").Value = "Boedo"
$ws.Range("N16").Value = "Capital Sur"
